# Adds the four new rule/antibiotic rows (Mastitis, Inpatient, Clean
# elective procedures, Contaminated procedures) to the info_texts sheet,
# matching the target commit "added all rules/antibiotics (except eye)".

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New shared strings are appended in first-use order, so these are entered
# in the same sequence the original author typed them (A5, B5, B6, A6, A7,
# A8, B7, B8) to keep xl/sharedStrings.xml's <si> ordering identical.

# Row 5: Mastitis
$ws.Range("A5").Value = "Mastitis"
$ws.Range("B5").Value = "Frequent combinations of antibiotics: Penicillin & Gentamicin (systemic & intramammary), Trimethoprim / Sulfonamid systemic or Penicillin / Neomycin / Trimethoprim / Sulfonamid. Udder injectors can be repurposed from farm animal medicine."

# Row 6: Inpatient
$ws.Range("B6").Value = "Caution: inpatient monitoring recommended."
$ws.Range("A6").Value = "Inpatient"

# Row 7 & 8 labels (column A), then their advice text (column B)
$ws.Range("A7").Value = "Clean elective procedures"
$ws.Range("A8").Value = "Contaminated procedures"
$ws.Range("B7").Value = "Recommended duration of antibiotics for clean, short elective surgical procedures with more than one person present is 24 hours."
$ws.Range("B8").Value = "Recommended duration of antibiotics for contaminated surgical procedures with more than one person present is 5 days."

# Match the row heights Excel computed for the wrapped text in each row.
$ws.Rows.Item(5).RowHeight = 85
$ws.Rows.Item(6).RowHeight = 17
$ws.Rows.Item(7).RowHeight = 51
$ws.Rows.Item(8).RowHeight = 51

# A latent, unreferenced conditional-format fill (same green swatch already
# lingering in styles.xml from earlier edits) gets duplicated as a dxf entry
# whenever the sheet is touched again in Excel; reproduce that by adding and
# immediately removing a matching highlight rule.
$fc = $ws.Range("A5").FormatConditions.Add(1, 3, "1")
$fc.Interior.Color = 13492663
$fc.Delete()

# The author's final action was entering the Contaminated-procedures advice
# in B8, leaving that cell selected.
$ws.Range("B8").Select()
